$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T55")

# Row 2 (Q0)
$ws.Range("B2").Value = 0.5271936044704619
$ws.Range("C2").Value = 0.5352467862865621
$ws.Range("D2").Value = 0.3759612643965448
$ws.Range("E2").Value = 0.6131568024547593
$ws.Range("F2").Value = 0.3249135209474168
$ws.Range("G2").Value = 14

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3608722674935116
$ws.Range("C3").Value = 0.4023274030500674
$ws.Range("D3").Value = 0.2193608715941454
$ws.Range("E3").Value = 0.4683597672667298
$ws.Range("F3").Value = 0.310740649621604
$ws.Range("G3").Value = 13

# Row 4 (Q2)
$ws.Range("B4").Value = 0.2850699742386379
$ws.Range("C4").Value = 0.3326698397266342
$ws.Range("D4").Value = 0.1632228408311082
$ws.Range("E4").Value = 0.4040084662864236
$ws.Range("F4").Value = 0.2990128315009373
$ws.Range("G4").Value = 12

# Row 5 (Q3)
$ws.Range("B5").Value = 0.437101989811909
$ws.Range("C5").Value = 0.4532318499641988
$ws.Range("D5").Value = 0.2559434172938412
$ws.Range("E5").Value = 0.5059085068407539
$ws.Range("F5").Value = 0.2671587441502564
$ws.Range("G5").Value = 11

# Row 6 (Q4)
$ws.Range("B6").Value = 0.3896036158490621
$ws.Range("C6").Value = 0.4137899396662769
$ws.Range("D6").Value = 0.213010691573719
$ws.Range("E6").Value = 0.4615308132440552
$ws.Range("F6").Value = 0.2608100928752895
$ws.Range("G6").Value = 10

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3372085974095211
$ws.Range("C7").Value = 0.3720251510852938
$ws.Range("D7").Value = 0.1694528389792911
$ws.Range("E7").Value = 0.4116464975914299
$ws.Range("F7").Value = 0.2504218459199277
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = 0.3687366187234252
$ws.Range("C8").Value = 0.3981931578693139
$ws.Range("D8").Value = 0.1958236157540623
$ws.Range("E8").Value = 0.4425196218859253
$ws.Range("F8").Value = 0.2615490596241174
$ws.Range("G8").Value = 8

# Row 9 (Q7)
$ws.Range("B9").Value = 0.3640834566367565
$ws.Range("C9").Value = 0.3961645473499545
$ws.Range("D9").Value = 0.1981723784143557
$ws.Range("E9").Value = 0.4451655629250265
$ws.Range("F9").Value = 0.2766795092775837
$ws.Range("G9").Value = 7

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3225455632229594
$ws.Range("C10").Value = 0.3512104856557821
$ws.Range("D10").Value = 0.1570729010663444
$ws.Range("E10").Value = 0.3963242372935882
$ws.Range("F10").Value = 0.2522790376821546
$ws.Range("G10").Value = 6

# Row 11 (Q9)
$ws.Range("B11").Value = 0.3733140426150657
$ws.Range("C11").Value = 0.3962272652827047
$ws.Range("D11").Value = 0.198344628268327
$ws.Range("E11").Value = 0.4453589880852603
$ws.Range("F11").Value = 0.2715263657886741
$ws.Range("G11").Value = 5
